# Sample-Exam.xlsx fix: attachment / playback-time format changes
# (commit "[DONE] Fix upload failed on CreateExercise.js (#65)")
#
#  - C3  (Exam Name/Description "Ujian Akhir Semester 2020/2021") gains a
#    trailing space.
#  - C9/C10 (question text) get "QWER" appended.
#  - D9/D10 (media link) values get wrapped as a single-element JSON array
#    string, e.g. ["https://...mp3"], instead of the bare URL - the
#    underlying hyperlink (still pointing at the bare URL) is left as-is.
#  - Active selection moves to D15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- C9 / C10: append "QWER" to the existing question text -------------
$c9 = $ws.Range("C9").Value()
$ws.Range("C9").Value = $c9 + "QWER"

$c10 = $ws.Range("C10").Value()
$ws.Range("C10").Value = $c10 + "QWER"

# --- D9 / D10: wrap the media URL in a JSON array ------------------------
$mediaUrl = "https://mcdn.podbean.com/mf/web/24psp6/2_Navigating_the_path_to_fluency64292.mp3"
$ws.Range("D9").Value = '["' + $mediaUrl + '"]'
$ws.Range("D10").Value = '["' + $mediaUrl + '"]'

# --- C3: "Ujian Akhir Semester 2020/2021" -> same text + trailing space ---
$ws.Range("C3").Value = "Ujian Akhir Semester 2020/2021 "

# --- Update the last-active-cell selection to match the saved state -----
$ws.Range("D15").Select()
